# "fix: end of day backup"
#
# The "Branching" test block (rows 24-36 of the "Train" section) gets two
# new rows inserted, a couple of the surrounding rows get their wording
# tweaked, and everything below (the "Entities" section) simply shifts
# down by two rows as a result. We reproduce that by inserting two blank
# rows at the right spots and then writing the new/changed cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the two new rows -------------------------------------------------
# New row 26: "Close"/"Delete" buttons change to "Save Branch"/"Abandon Branch"
$ws.Rows("26:26").Insert()
# New row 32: "After branch is abandonded the original training should remain
# unchanged" (row numbers below this point are already shifted down by the
# first insert above)
$ws.Rows("32:32").Insert()

# --- 2. Fill in the newly inserted / reworded cells -----------------------------

# Row 26 (new) - mirrors the Train / Branching / Train Dialog pattern of its neighbours
$ws.Range("A26").Value = "Train"
$ws.Range("B26").Value = "Branching"
$ws.Range("C26").Value = "Train Dialog"
$ws.Range("D26").Value = '"Close" and "Delete" buttons should change to "Save Branch" and "Abandon Branch" after branching'

# Row 31 (formerly row 30, shifted by the first insert) - wording changed, and a
# new "manual" note added in column E
$ws.Range("D31").Value = "After branch is saved the original training should remain unchanged"
$ws.Range("E31").Value = "manual"

# Row 32 (new)
$ws.Range("A32").Value = "Train"
$ws.Range("B32").Value = "Branching"
$ws.Range("C32").Value = "Train Dialog"
$ws.Range("D32").Value = "After branch is abandonded the original training should remain unchanged"
$ws.Range("E32").Value = "manual"

# --- 3. Resize the table / autofilter to cover the two new rows ----------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F59"))

# --- 4. Restore the selection state (cosmetic) ----------------------------------
$ws.Range("D33").Select()
